# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-09-12 (serial 45181) to 2023-09-13 (serial 45182).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 331
$col = 3  # Column C

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    if ($cell.Value2 -eq 45181) {
        $cell.Value2 = 45182
    }
}
